$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before current column D ("Tipo") to make room for "MAE"
$ws.Range("D1").EntireColumn.Insert()

# New header cell D1 = "MAE", matching the style of the other header cells
$ws.Range("D1").Value = "MAE"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1

# Update values per diff
$ws.Range("B2").Value = 0.05688602116888641
$ws.Range("C2").Value = 0.9992297983880002
$ws.Range("D2").Value = 0.1912640313995153
$ws.Range("E2").Value = "single"
